$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("A1").Value = "lxmpqrh@gzcmhid.com"
$ws.Range("B1").Value = "zp89Pr2"
$ws.Range("C1").Value = "gefkqyw"
$ws.Range("D1").Value = "Travell"

# Row 2
$ws.Range("A2").Value = "uytvxnsk@aheblyml.com"
$ws.Range("B2").Value = "IO7vYRsw"
$ws.Range("C2").Value = "lcldcrhi"
$ws.Range("D2").Value = "Software"

# Row 3
$ws.Range("A3").Value = "rguvqduo@lssjcxvl.com"
$ws.Range("B3").Value = "5b3RISBIY7N"
$ws.Range("C3").Value = "tctiioxxlgm"
$ws.Range("D3").Value = "I love to code in Py"

# Row 4
$ws.Range("A4").Value = "dhkflsmx@mdksplfs.com"
$ws.Range("B4").Value = "wag6JQNc"
$ws.Range("C4").Value = "pdguxjge"
$ws.Range("D4").Value = "Foodie. "

# Row 5
$ws.Range("A5").Value = "zhbggvvn@ybnxlseg.com"
$ws.Range("B5").Value = "YAceJZ"
$ws.Range("C5").Value = "tipuft"
$ws.Range("D5").Value = "Foodie. Yoga enthusi"
